$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 11.29255533333333
$ws.Range("H2").Value = 33.877666
$ws.Range("I2").Value = 0.5495662219753726
$ws.Range("J2").Value = 0.6375557499803809
$ws.Range("M2").Value = 29.785352
$ws.Range("N2").Value = 89.356056
$ws.Range("O2").Value = 0.7923195065866085
$ws.Range("P2").Value = 0.7947519366640845
$ws.Range("Q2").Value = 336.3527355828107
$ws.Range("R2").Value = 3027.174620245296
$ws.Range("S2").Value = 0.4354320378321938
$ws.Range("T2").Value = 0.5066986670282305
$ws.Range("G3").Value = 11.29255533333333
$ws.Range("H3").Value = 33.877666
$ws.Range("I3").Value = 0.5495662219753726
$ws.Range("J3").Value = 0.6375557499803809
$ws.Range("O3").Value = 0.1171985110386058
$ws.Range("P3").Value = 0.1175583118271966
$ws.Range("Q3").Value = 49.752706409934
$ws.Range("R3").Value = 447.774357689406
$ws.Range("S3").Value = 0.06440834293262561
$ws.Range("T3").Value = 0.0749499776634158
$ws.Range("G4").Value = 11.29255533333333
$ws.Range("H4").Value = 33.877666
$ws.Range("I4").Value = 0.5495662219753726
$ws.Range("J4").Value = 0.6375557499803809
$ws.Range("M4").Value = 1.362560333333333
$ws.Range("N4").Value = 4.087681
$ws.Range("O4").Value = 0.03624543805965938
$ws.Range("P4").Value = 0.03635671197501131
$ws.Range("Q4").Value = 15.38678795917178
$ws.Range("R4").Value = 138.481091632546
$ws.Range("S4").Value = 0.01991926845828939
$ws.Range("T4").Value = 0.02317943077004903
$ws.Range("G5").Value = 11.29255533333333
$ws.Range("H5").Value = 33.877666
$ws.Range("I5").Value = 0.5495662219753726
$ws.Range("J5").Value = 0.6375557499803809
$ws.Range("M5").Value = 0.3451695
$ws.Range("N5").Value = 0.690339
$ws.Range("O5").Value = 0.00918184643004207
$ws.Range("P5").Value = 0.006140023203404898
$ws.Range("Q5").Value = 3.897845678129001
$ws.Range("R5").Value = 23.38707406877401
$ws.Range("S5").Value = 0.005046032653316283
$ws.Range("T5").Value = 0.003914607098343751
$ws.Range("G6").Value = 11.29255533333333
$ws.Range("H6").Value = 33.877666
$ws.Range("I6").Value = 0.5495662219753726
$ws.Range("J6").Value = 0.6375557499803809
$ws.Range("M6").Value = 1.693723333333333
$ws.Range("N6").Value = 5.08117
$ws.Range("O6").Value = 0.04505469788508434
$ws.Range("P6").Value = 0.04519301633030275
$ws.Range("Q6").Value = 19.12646446102445
$ws.Range("R6").Value = 172.13818014922
$ws.Range("S6").Value = 0.02476054009894762
$ws.Range("T6").Value = 0.02881306742034177
$ws.Range("I7").Value = 0.03478077306145753
$ws.Range("J7").Value = 0.04034942645199305
$ws.Range("M7").Value = 29.785352
$ws.Range("N7").Value = 89.356056
$ws.Range("O7").Value = 0.7923195065866085
$ws.Range("P7").Value = 0.7947519366640845
$ws.Range("Q7").Value = 21.28698543890933
$ws.Range("R7").Value = 191.582868950184
$ws.Range("S7").Value = 0.02755748495075484
$ws.Range("T7").Value = 0.03206778481600652
$ws.Range("I8").Value = 0.03478077306145753
$ws.Range("J8").Value = 0.04034942645199305
$ws.Range("O8").Value = 0.1171985110386058
$ws.Range("P8").Value = 0.1175583118271966
$ws.Range("S8").Value = 0.004076254815574475
$ws.Range("T8").Value = 0.004743410456891934
$ws.Range("I9").Value = 0.03478077306145753
$ws.Range("J9").Value = 0.04034942645199305
$ws.Range("M9").Value = 1.362560333333333
$ws.Range("N9").Value = 4.087681
$ws.Range("O9").Value = 0.03624543805965938
$ws.Range("P9").Value = 0.03635671197501131
$ws.Range("Q9").Value = 0.973794164839889
$ws.Range("R9").Value = 8.764147483559
$ws.Range("S9").Value = 0.001260644355666128
$ws.Range("T9").Value = 0.001466972475872014
$ws.Range("I10").Value = 0.03478077306145753
$ws.Range("J10").Value = 0.04034942645199305
$ws.Range("M10").Value = 0.3451695
$ws.Range("N10").Value = 0.690339
$ws.Range("O10").Value = 0.00918184643004207
$ws.Range("P10").Value = 0.006140023203404898
$ws.Range("Q10").Value = 0.2466856232035
$ws.Range("R10").Value = 1.480113739221
$ws.Range("S10").Value = 0.0003193517169684473
$ws.Range("T10").Value = 0.0002477464146593167
$ws.Range("I11").Value = 0.03478077306145753
$ws.Range("J11").Value = 0.04034942645199305
$ws.Range("M11").Value = 1.693723333333333
$ws.Range("N11").Value = 5.08117
$ws.Range("O11").Value = 0.04505469788508434
$ws.Range("P11").Value = 0.04519301633030275
$ws.Range("Q11").Value = 1.210469627292222
$ws.Range("R11").Value = 10.89422664563
$ws.Range("S11").Value = 0.001567037222493649
$ws.Range("T11").Value = 0.001823512288563272
$ws.Range("G12").Value = 8.5075845
$ws.Range("H12").Value = 17.015169
$ws.Range("I12").Value = 0.4140321595768645
$ws.Range("J12").Value = 0.3202144691088791
$ws.Range("M12").Value = 29.785352
$ws.Range("N12").Value = 89.356056
$ws.Range("O12").Value = 0.7923195065866085
$ws.Range("P12").Value = 0.7947519366640845
$ws.Range("Q12").Value = 253.401399002244
$ws.Range("R12").Value = 1520.408394013464
$ws.Range("S12").Value = 0.3280457563869293
$ws.Range("T12").Value = 0.2544910694721433
$ws.Range("G13").Value = 8.5075845
$ws.Range("H13").Value = 17.015169
$ws.Range("I13").Value = 0.4140321595768645
$ws.Range("J13").Value = 0.3202144691088791
$ws.Range("O13").Value = 0.1171985110386058
$ws.Range("P13").Value = 0.1175583118271966
$ws.Range("Q13").Value = 37.4826902673465
$ws.Range("R13").Value = 224.896141604079
$ws.Range("S13").Value = 0.04852395262450697
$ws.Range("T13").Value = 0.03764387241108182
$ws.Range("G14").Value = 8.5075845
$ws.Range("H14").Value = 17.015169
$ws.Range("I14").Value = 0.4140321595768645
$ws.Range("J14").Value = 0.3202144691088791
$ws.Range("M14").Value = 1.362560333333333
$ws.Range("N14").Value = 4.087681
$ws.Range("O14").Value = 0.03624543805965938
$ws.Range("P14").Value = 0.03635671197501131
$ws.Range("Q14").Value = 11.5920971721815
$ws.Range("R14").Value = 69.55258303308899
$ws.Range("S14").Value = 0.01500677699465025
$ws.Range("T14").Value = 0.01164194522362268
$ws.Range("G15").Value = 8.5075845
$ws.Range("H15").Value = 17.015169
$ws.Range("I15").Value = 0.4140321595768645
$ws.Range("J15").Value = 0.3202144691088791
$ws.Range("M15").Value = 0.3451695
$ws.Range("N15").Value = 0.690339
$ws.Range("O15").Value = 0.00918184643004207
$ws.Range("P15").Value = 0.006140023203404898
$ws.Range("Q15").Value = 2.93655868807275
$ws.Range("R15").Value = 11.746234752291
$ws.Range("S15").Value = 0.003801579706333442
$ws.Range("T15").Value = 0.001966124270394499
$ws.Range("G16").Value = 8.5075845
$ws.Range("H16").Value = 17.015169
$ws.Range("I16").Value = 0.4140321595768645
$ws.Range("J16").Value = 0.3202144691088791
$ws.Range("M16").Value = 1.693723333333333
$ws.Range("N16").Value = 5.08117
$ws.Range("O16").Value = 0.04505469788508434
$ws.Range("P16").Value = 0.04519301633030275
$ws.Range("Q16").Value = 14.409494377955
$ws.Range("R16").Value = 86.45696626773001
$ws.Range("S16").Value = 0.01865409386444466
$ws.Range("T16").Value = 0.0144714577316368
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.03330533333333333
$ws.Range("H17").Value = 0.099916
$ws.Range("I17").Value = 0.001620845386305282
$ws.Range("J17").Value = 0.001880354458746944
$ws.Range("M17").Value = 29.785352
$ws.Range("N17").Value = 89.356056
$ws.Range("O17").Value = 0.7923195065866085
$ws.Range("P17").Value = 0.7947519366640845
$ws.Range("Q17").Value = 0.9920110768106666
$ws.Range("R17").Value = 8.928099691296
$ws.Range("S17").Value = 0.001284227416730582
$ws.Range("T17").Value = 0.00149441534770408
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.03330533333333333
$ws.Range("H18").Value = 0.099916
$ws.Range("I18").Value = 0.001620845386305282
$ws.Range("J18").Value = 0.001880354458746944
$ws.Range("O18").Value = 0.1171985110386058
$ws.Range("P18").Value = 0.1175583118271966
$ws.Range("Q18").Value = 0.146736537684
$ws.Range("R18").Value = 1.320628839156
$ws.Range("S18").Value = 0.0001899606658987729
$ws.Range("T18").Value = 0.0002210512958070326
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.03330533333333333
$ws.Range("H19").Value = 0.099916
$ws.Range("I19").Value = 0.001620845386305282
$ws.Range("J19").Value = 0.001880354458746944
$ws.Range("M19").Value = 1.362560333333333
$ws.Range("N19").Value = 4.087681
$ws.Range("O19").Value = 0.03624543805965938
$ws.Range("P19").Value = 0.03635671197501131
$ws.Range("Q19").Value = 0.04538052608844444
$ws.Range("R19").Value = 0.408424734796
$ws.Range("S19").Value = 0.00005874825105361278
$ws.Range("T19").Value = 0.00006836350546759092
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.03330533333333333
$ws.Range("H20").Value = 0.099916
$ws.Range("I20").Value = 0.001620845386305282
$ws.Range("J20").Value = 0.001880354458746944
$ws.Range("M20").Value = 0.3451695
$ws.Range("N20").Value = 0.690339
$ws.Range("O20").Value = 0.00918184643004207
$ws.Range("P20").Value = 0.006140023203404898
$ws.Range("Q20").Value = 0.011495985254
$ws.Range("R20").Value = 0.06897591152400001
$ws.Range("S20").Value = 0.00001488235342389731
$ws.Range("T20").Value = 0.00001154542000733209
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.03330533333333333
$ws.Range("H21").Value = 0.099916
$ws.Range("I21").Value = 0.001620845386305282
$ws.Range("J21").Value = 0.001880354458746944
$ws.Range("M21").Value = 1.693723333333333
$ws.Range("N21").Value = 5.08117
$ws.Range("O21").Value = 0.04505469788508434
$ws.Range("P21").Value = 0.04519301633030275
$ws.Range("Q21").Value = 0.05641002019111111
$ws.Range("R21").Value = 0.50769018172
$ws.Range("S21").Value = 0.00007302669919841731
$ws.Range("T21").Value = 0.0002817888967396490
